$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.185198527466267
$ws.Range("C2").Value = 0.6817654794757074
$ws.Range("D2").Value = 0.8357345580027106
$ws.Range("E2").Value = 0.9141851880241282
$ws.Range("F2").Value = 0.9153496465291245
$ws.Range("G2").Value = 23

$ws.Range("B3").Value = 0.6061131440253159
$ws.Range("C3").Value = 0.9235219194524696
$ws.Range("D3").Value = 2.046785243796049
$ws.Range("E3").Value = 1.430659024294765
$ws.Range("F3").Value = 1.32641777175159
$ws.Range("G3").Value = 22

$ws.Range("B4").Value = 0.4736973475431321
$ws.Range("C4").Value = 1.172664837566667
$ws.Range("D4").Value = 3.405130736559959
$ws.Range("E4").Value = 1.845299633273675
$ws.Range("F4").Value = 1.827506125151183
$ws.Range("G4").Value = 21

$ws.Range("B5").Value = 0.5953765531118547
$ws.Range("C5").Value = 0.8086914423054233
$ws.Range("D5").Value = 1.001718343406927
$ws.Range("E5").Value = 1.000858802932225
$ws.Range("F5").Value = 0.8254154318705684
$ws.Range("G5").Value = 20

$ws.Range("B6").Value = 0.4684946891824553
$ws.Range("C6").Value = 0.6936936078512286
$ws.Range("D6").Value = 0.7365785875710174
$ws.Range("E6").Value = 0.8582415671423852
$ws.Range("F6").Value = 0.7387953769405898
$ws.Range("G6").Value = 19

$ws.Range("B7").Value = 0.2823731533649205
$ws.Range("C7").Value = 0.5667961785861276
$ws.Range("D7").Value = 0.457481639405574
$ws.Range("E7").Value = 0.6763738902453095
$ws.Range("F7").Value = 0.6324298031243258
$ws.Range("G7").Value = 18

$ws.Range("B8").Value = 0.2744592162286795
$ws.Range("C8").Value = 0.5038470190539809
$ws.Range("D8").Value = 0.354955400894988
$ws.Range("E8").Value = 0.5957813364775604
$ws.Range("F8").Value = 0.5450727114268883
$ws.Range("G8").Value = 17

$ws.Range("B9").Value = 0.2856656093728968
$ws.Range("C9").Value = 0.4249183470701322
$ws.Range("D9").Value = 0.2447576538564426
$ws.Range("E9").Value = 0.4947298796883433
$ws.Range("F9").Value = 0.4218825517020163
$ws.Range("G9").Value = 12

$ws.Range("B10").Value = 0.1512223912673487
$ws.Range("C10").Value = 0.384830859019198
$ws.Range("D10").Value = 0.2185051035033979
$ws.Range("E10").Value = 0.467445294663876
$ws.Range("F10").Value = 0.4777478838570751
$ws.Range("G10").Value = 7

$ws.Range("B11").Value = 0.1070798832976198
$ws.Range("C11").Value = 0.5676804838982221
$ws.Range("D11").Value = 0.3539232406442811
$ws.Range("E11").Value = 0.5949144817906865
$ws.Range("F11").Value = 0.7167187097152367
